$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 71 ("Vega Modelo de Temuco" / Achicoria
# weekly report gains a new latest-week entry). This pushes the existing
# rows 71 and 72 down to 72 and 73 respectively, preserving their data.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with this week's data, matching the
# layout of the surrounding rows.
$ws.Cells.Item(71, 1).Value = 10
$ws.Cells.Item(71, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(71, 3).Value = "La Araucanía"
$ws.Cells.Item(71, 4).Value = 45013
$ws.Cells.Item(71, 4).NumberFormat = $ws.Cells.Item(72, 4).NumberFormat
$ws.Cells.Item(71, 5).Value = 9
$ws.Cells.Item(71, 6).Value = 100112010
$ws.Cells.Item(71, 7).Value = "Achicoria"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 50
$ws.Cells.Item(71, 11).Value = 10000
$ws.Cells.Item(71, 12).Value = 10000
$ws.Cells.Item(71, 13).Value = 10000
$ws.Cells.Item(71, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(71, 15).Value = "Región Metropolitana"
$ws.Cells.Item(71, 16).Value = 556
$ws.Cells.Item(71, 17).Value = 18
$ws.Cells.Item(71, 18).Value = "Hortaliza"
